$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fix the duplicate/placeholder score value (-1) to the real numeric score
$ws.Range("B2").Value = 1000

# Row 3: add a new ID/score entry
$ws.Range("A3").Value = 1397200151
$ws.Range("B3").Value = "'6"
$ws.Range("B3").Style = "Normal"

# Leave selection on the newly added cell, like the author did
[void]$ws.Range("B3").Select()
